$wb = $excel.ActiveWorkbook

# Add a new worksheet at the very end of the workbook (after the last sheet)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "rty"

# Fill in the header row
$newSheet.Range("A1").Value = "Outdoor Model"
$newSheet.Range("B1").Value = "Outdoor Quantity"
$newSheet.Range("C1").Value = "Outdoor Serial(s)"
$newSheet.Range("D1").Value = "Indoor Model"
$newSheet.Range("E1").Value = "Indoor Quantity"
$newSheet.Range("F1").Value = "Indoor Serial(s)"

# Apply bold + centered + bordered style like other sheet headers
$headerRange = $newSheet.Range("A1:F1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1
